# Sprint 0 presentatie - puntjes
#
# Slide 6 ("Wat hebben we gedaan?"), the bullet that reads "Product backlog"
# becomes "Concept product backlog" - only the first run of that paragraph
# ("Product ") is touched; the second run ("backlog", which carries the
# spell-check err="1" flag) is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(3)
$run = $para.Runs(1)
$run.Text = "Concept product "
